$wb = $excel.ActiveWorkbook

# Add the new "Tasks" worksheet after the last existing sheet ("Deals")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tasks"

# Header row
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Completion"

# Data rows
$ws.Range("A2").Value = "BVC "
$ws.Range("B2").Value = 25

$ws.Range("A3").Value = "ATCO Lab"
$ws.Range("B3").Value = 65

$ws.Range("A4").Value = "General Dynamics"
$ws.Range("B4").Value = 10

# Column widths (engine stores width as ColumnWidth + 5/6; pre-compensate so the
# serialized <col width> lands on the recorded values)
$ws.Columns.Item(1).ColumnWidth = 17.498697916666668
$ws.Columns.Item(2).ColumnWidth = 15.944010416666666

# Selection matches the recorded view state
$ws.Range("D3").Select() | Out-Null
